# Fruta / hortaliza, semanal
#
# Two new daily price records were added to the "Vega Monumental
# Concepción - Piña" sheet. Inserting them in the middle of the existing
# (date-unsorted) list pushes every following record down by one row
# each time, which is why the unified diff shows almost every row from
# 77 onward "changing" - in reality nothing about those pre-existing
# records changed, they just moved down.
#
# Net effect:
#   * insert one new row at row 77 (previous rows 77..129 become 78..130)
#   * insert a second new row at (the now shifted) row 123
#     (previous rows 123..130 become 124..131)
#   * the sheet's used range grows from A1:T129 to A1:T131
#
# Columns A, B, C, E, F, G, H, I, J, R are constant across the whole
# sheet (same market / product taxonomy for every row), so the two new
# rows reuse those same values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert first new row at 77 -------------------------------------
$ws.Rows.Item(77).Insert()

$ws.Range("A77").Value = 11
$ws.Range("B77").Value = "Vega Monumental Concepción"
$ws.Range("C77").Value = "Bíobío"
$ws.Range("D77").Value = 44567
$ws.Range("E77").Value = 8
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108005
$ws.Range("J77").Value = "Piña"
$ws.Range("K77").Value = "Caramelo"
$ws.Range("L77").Value = "Segunda"
$ws.Range("M77").Value = 310
$ws.Range("N77").Value = 14000
$ws.Range("O77").Value = 15000
$ws.Range("P77").Value = 14484
$ws.Range("Q77").Value = "$/caja 14 unidades"
$ws.Range("R77").Value = "Ecuador"
$ws.Range("S77").Value = 1035
$ws.Range("T77").Value = 14

# --- insert second new row at (shifted) 123 --------------------------
$ws.Rows.Item(123).Insert()

$ws.Range("A123").Value = 11
$ws.Range("B123").Value = "Vega Monumental Concepción"
$ws.Range("C123").Value = "Bíobío"
$ws.Range("D123").Value = 44568
$ws.Range("E123").Value = 8
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = "Tropicales y subtropicales"
$ws.Range("I123").Value = 100108005
$ws.Range("J123").Value = "Piña"
$ws.Range("K123").Value = "Sin especificar"
$ws.Range("L123").Value = "Segunda"
$ws.Range("M123").Value = 290
$ws.Range("N123").Value = 15000
$ws.Range("O123").Value = 16000
$ws.Range("P123").Value = 15483
$ws.Range("Q123").Value = "$/caja 14 unidades"
$ws.Range("R123").Value = "Ecuador"
$ws.Range("S123").Value = 1106
$ws.Range("T123").Value = 14
